# Witze_Styger.docx edit:
#   - The "SW7:" joke block stays as-is, but loses the (hidden) _GoBack
#     bookmark that used to sit on its heading paragraph.
#   - A brand-new "SW8:" joke block is appended at the very end of the
#     document (blank separator line, bold "SW8:" heading, the two joke
#     lines, and the _GoBack bookmark now sitting at the end of the very
#     last sentence) followed by a trailing blank paragraph.
#
# We locate the "SW7:" heading paragraph and replace everything from its
# start through the end of the document with the fully-specified target
# markup (InsertXML replaces the contents of the range it is called on).

$d = $word.ActiveDocument

$sw7 = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    $t = $t.Replace([char]13, "").Replace([char]7, "")
    if ($t -eq "SW7:") {
        $sw7 = $p
    }
}

$startPos = $sw7.Range.Start
$endPos = $d.Content.End
$target = $d.Range($startPos, $endPos)

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val="de-CH"/></w:rPr><w:lastRenderedPageBreak/><w:t>SW7:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Tierhandlung;</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>C-Affe = 5000 Fr</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>C++Affe = 10000 Fr.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Affe ohne Schild = 50000 Fr.,  macht nie etwas = Berater</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:b/><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val="de-CH"/></w:rPr><w:t>SW8:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve">2 Ing am </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Bhf</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Freundin will zurück zur Mama, Zettel am Kühlschrank</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr></w:p>'

$target.InsertXML($newXml)
